$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "clean" number-looking string need to be
# forced to Text format first, otherwise Excel auto-converts them to a
# real number (losing e.g. trailing zeros: "0.1000" -> 0.1).
$textCells = @(
    "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D18", "D21", "D22", "D23", "D24", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D46", "D47", "D48", "D49", "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.402.36'
$ws.Range("E2").Value = '  -2.81%  '

$ws.Range("D3").Value = '2.274.22'
$ws.Range("E3").Value = '  -4.70%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '299.96'
$ws.Range("E5").Value = '  -3.27%  '

$ws.Range("D6").Value = '96.80'
$ws.Range("E6").Value = '  -7.46%  '

$ws.Range("D7").Value = '0.504'
$ws.Range("E7").Value = '  -1.27%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = '0.497'
$ws.Range("E9").Value = '  -4.57%  '

$ws.Range("D10").Value = '33.82'
$ws.Range("E10").Value = '  -6.32%  '

$ws.Range("D11").Value = '50.65'
$ws.Range("E11").Value = '  -5.15%  '

$ws.Range("D12").Value = '0.0786'
$ws.Range("E12").Value = '  -3.45%  '

$ws.Range("E13").Value = '  +0.04%  '

$ws.Range("E14").Value = '  -5.07%  '

$ws.Range("D15").Value = '2.624.70'
$ws.Range("E15").Value = '  -4.73%  '

$ws.Range("E16").Value = '  -3.01%  '

$ws.Range("D17").Value = '2.273.57'
$ws.Range("E17").Value = '  -4.54%  '

$ws.Range("D18").Value = '0.784'
$ws.Range("E18").Value = '  -3.43%  '

$ws.Range("D19").Value = '42.304.65'
$ws.Range("E19").Value = '  -2.99%  '

$ws.Range("D20").Value = '0.0₃0893'
$ws.Range("E20").Value = '  -2.68%  '

$ws.Range("D21").Value = '11.43'
$ws.Range("E21").Value = '  -4.04%  '

$ws.Range("D22").Value = '5.99'
$ws.Range("E22").Value = '  -5.39%  '

$ws.Range("D23").Value = '66.56'
$ws.Range("E23").Value = '  -2.77%  '

$ws.Range("D24").Value = '234.14'
$ws.Range("E24").Value = '  -2.89%  '

$ws.Range("E25").Value = '  -5.99%  '

$ws.Range("E26").Value = '  -5.31%  '

$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("D28").Value = '24.36'
$ws.Range("E28").Value = '  -5.93%  '

$ws.Range("D29").Value = '2.07'
$ws.Range("E29").Value = '  -1.73%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '33.83'
$ws.Range("E30").Value = '  -7.78%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '163.79'
$ws.Range("E31").Value = '  +1.84%  '

$ws.Range("D32").Value = '9.07'
$ws.Range("E32").Value = '  -4.88%  '

$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("D34").Value = '4.94'
$ws.Range("E34").Value = '  -6.06%  '

$ws.Range("E35").Value = '  -4.83%  '

$ws.Range("D36").Value = '0.0695'
$ws.Range("E36").Value = '  -6.15%  '

$ws.Range("D37").Value = '4.37'
$ws.Range("E37").Value = '  -6.47%  '

$ws.Range("D38").Value = '2.82'
$ws.Range("E38").Value = '  -9.92%  '

$ws.Range("D39").Value = '16.08'
$ws.Range("E39").Value = '  -12.55%  '

$ws.Range("D40").Value = '0.1000'
$ws.Range("E40").Value = '  -5.63%  '

$ws.Range("E41").Value = '  -9.39%  '

$ws.Range("D42").Value = '0.110'
$ws.Range("E42").Value = '  -3.74%  '

$ws.Range("D43").Value = '2.39'
$ws.Range("E43").Value = '  -7.45%  '

$ws.Range("D44").Value = '1.964.62'
$ws.Range("E44").Value = '  -3.57%  '

$ws.Range("E45").Value = '  -3.35%  '

$ws.Range("D46").Value = '17.92'
$ws.Range("E46").Value = '  -9.45%  '

$ws.Range("D47").Value = '9.65'
$ws.Range("E47").Value = '  -8.76%  '

$ws.Range("D48").Value = '2.82'
$ws.Range("E48").Value = '  -10.01%  '

$ws.Range("D49").Value = '2.83'
$ws.Range("E49").Value = '  -4.59%  '

$ws.Range("D50").Value = '4.68'
$ws.Range("E50").Value = '  -1.40%  '

$ws.Range("D51").Value = '2.499.45'
$ws.Range("E51").Value = '  -4.41%  '
